$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalidCredentialTest")
$ws.Activate()

# Change A2 from the text "Peter" to the number 1122
$ws.Range("A2").Value = 1122

# Update the selection to reflect the new active cell
$ws.Range("A2").Select()
